$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume cells hold text that can look numeric (e.g. "44.35",
# "96.958.68"). Force each touched cell to Text format right before writing
# so Excel does not silently reinterpret it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.958.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.677.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +20.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '657.43'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.424'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.72%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.678.42'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.35'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.205'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.52'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.356.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.644.01'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.664.18'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.78'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.40'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.531'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.98%  '
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.46'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '513.44'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.43'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.10'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.71%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +14.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.04'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.96'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.75%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.590'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '617.76'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '42.74'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +25.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.160'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.963'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.60%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.73%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.16'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +8.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0441'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.420'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +26.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.31'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.62'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.61'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.94%  '
